# AIC9parameters.xlsx was regenerated (e.g. re-run of the MATLAB script that
# produces this data) so the 16x3 numeric grid on Sheet1 has new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2204.2292170912615
$ws.Range("B1").Value = 1383.0795960930136
$ws.Range("C1").Value = 1428.8271839685574

$ws.Range("A2").Value = 2227.1192368807488
$ws.Range("B2").Value = 1465.9570461612841
$ws.Range("C2").Value = 1494.2946379747702

$ws.Range("A3").Value = 2343.9137069530807
$ws.Range("B3").Value = 1586.5214891726905
$ws.Range("C3").Value = 1455.6163657279224

$ws.Range("A4").Value = 2320.962498724236
$ws.Range("B4").Value = 1774.6466847242061
$ws.Range("C4").Value = 1705.8701555345397

$ws.Range("A5").Value = 2422.6539933980139
$ws.Range("B5").Value = 1671.2731630940148
$ws.Range("C5").Value = 1620.8443032066889

$ws.Range("A6").Value = 2360.7084387163345
$ws.Range("B6").Value = 1774.6259027478818
$ws.Range("C6").Value = 1784.2650066614551

$ws.Range("A7").Value = 1992.98484726981
$ws.Range("B7").Value = 1566.3791362128925
$ws.Range("C7").Value = 1482.5142659959436

$ws.Range("A8").Value = 2135.7549744213029
$ws.Range("B8").Value = 1660.8003328907289
$ws.Range("C8").Value = 1524.1020273468509

$ws.Range("A9").Value = 2471.0153838041033
$ws.Range("B9").Value = 1788.3782491419802
$ws.Range("C9").Value = 1514.0343486951572

$ws.Range("A10").Value = 2111.5944041864032
$ws.Range("B10").Value = 1368.6689139450625
$ws.Range("C10").Value = 1292.6384240066725

$ws.Range("A11").Value = 1970.0038839175693
$ws.Range("B11").Value = 1416.2687623358067
$ws.Range("C11").Value = 1298.5132524254782

$ws.Range("A12").Value = 2787.9536320270217
$ws.Range("B12").Value = 2270.4512025612389
$ws.Range("C12").Value = 2036.3897228392641

$ws.Range("A13").Value = 2315.185341984487
$ws.Range("B13").Value = 1782.3538751813696
$ws.Range("C13").Value = 1792.0065844647258

$ws.Range("A14").Value = 2593.0484738491914
$ws.Range("B14").Value = 1922.7202893800636
$ws.Range("C14").Value = 1702.405037810182

$ws.Range("A15").Value = 2508.6641344546319
$ws.Range("B15").Value = 2033.5167339667075
$ws.Range("C15").Value = 1844.9470529650316

$ws.Range("A16").Value = 2205.5211166450599
$ws.Range("B16").Value = 1511.4651336237885
$ws.Range("C16").Value = 1268.1754357637358

# Cosmetic window-position metadata (bookViews/workbookView xWindow & yWindow
# 1152 -> 2304) also changed in the source diff. This is pure UI chrome state
# (the on-screen window position), not part of the Excel Range/Worksheet
# object model, so it is set here via the Application/Window position in
# case the host environment maps it through; it has no effect on workbook
# data either way.
try {
    $excel.Left = 2304 / 20
    $excel.Top = 2304 / 20
} catch {}
try {
    $win = $excel.ActiveWindow
    $win.Left = 2304 / 20
    $win.Top = 2304 / 20
} catch {}
